$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value2 = 3.748362666666667
$ws.Range("H2").Value2 = 11.245088
$ws.Range("I2").Value2 = 0.2072014058556041
$ws.Range("J2").Value2 = 0.2072014058556041
$ws.Range("K2").Value2 = 1
$ws.Range("L2").Value2 = 0.3333333333333333
$ws.Range("M2").Value2 = 0.06686833333333334
$ws.Range("N2").Value2 = 0.200605
$ws.Range("O2").Value2 = 0.05567274787007094
$ws.Range("P2").Value2 = 0.05567274787007093
$ws.Range("Q2").Value2 = 0.2506467642488889
$ws.Range("R2").Value2 = 2.25582087824
$ws.Range("S2").Value2 = 0.01153547162652329
$ws.Range("T2").Value2 = 0.01153547162652329
$ws.Range("G3").Value2 = 3.748362666666667
$ws.Range("H3").Value2 = 11.245088
$ws.Range("I3").Value2 = 0.2072014058556041
$ws.Range("J3").Value2 = 0.2072014058556041
$ws.Range("O3").Value2 = 0.8610905203551533
$ws.Range("P3").Value2 = 0.8610905203551533
$ws.Range("Q3").Value2 = 3.876754083633777
$ws.Range("R3").Value2 = 34.890786752704
$ws.Range("S3").Value2 = 0.1784191663865214
$ws.Range("T3").Value2 = 0.1784191663865215
$ws.Range("G4").Value2 = 3.748362666666667
$ws.Range("H4").Value2 = 11.245088
$ws.Range("I4").Value2 = 0.2072014058556041
$ws.Range("J4").Value2 = 0.2072014058556041
$ws.Range("O4").Value2 = 0.08323673177477579
$ws.Range("P4").Value2 = 0.08323673177477578
$ws.Range("Q4").Value2 = 0.3747438070542223
$ws.Range("R4").Value2 = 3.372694263488
$ws.Range("S4").Value2 = 0.01724676784255938
$ws.Range("T4").Value2 = 0.01724676784255938
$ws.Range("I5").Value2 = 0.6140431114114622
$ws.Range("J5").Value2 = 0.6140431114114623
$ws.Range("K5").Value2 = 1
$ws.Range("L5").Value2 = 0.3333333333333333
$ws.Range("M5").Value2 = 0.06686833333333334
$ws.Range("N5").Value2 = 0.200605
$ws.Range("O5").Value2 = 0.05567274787007094
$ws.Range("P5").Value2 = 0.05567274787007093
$ws.Range("Q5").Value2 = 0.7427937969294444
$ws.Range("R5").Value2 = 6.685144172364999
$ws.Range("S5").Value2 = 0.03418546732296422
$ws.Range("T5").Value2 = 0.03418546732296422
$ws.Range("I6").Value2 = 0.6140431114114622
$ws.Range("J6").Value2 = 0.6140431114114623
$ws.Range("O6").Value2 = 0.8610905203551533
$ws.Range("P6").Value2 = 0.8610905203551533
$ws.Range("S6").Value2 = 0.5287467023257933
$ws.Range("T6").Value2 = 0.5287467023257935
$ws.Range("I7").Value2 = 0.6140431114114622
$ws.Range("J7").Value2 = 0.6140431114114623
$ws.Range("O7").Value2 = 0.08323673177477579
$ws.Range("P7").Value2 = 0.08323673177477578
$ws.Range("S7").Value2 = 0.05111094176270465
$ws.Range("T7").Value2 = 0.05111094176270465
$ws.Range("G8").Value2 = 3.233763666666667
$ws.Range("H8").Value2 = 9.701291000000001
$ws.Range("I8").Value2 = 0.1787554827329337
$ws.Range("J8").Value2 = 0.1787554827329337
$ws.Range("K8").Value2 = 1
$ws.Range("L8").Value2 = 0.3333333333333333
$ws.Range("M8").Value2 = 0.06686833333333334
$ws.Range("N8").Value2 = 0.200605
$ws.Range("O8").Value2 = 0.05567274787007094
$ws.Range("P8").Value2 = 0.05567274787007093
$ws.Range("Q8").Value2 = 0.2162363867838889
$ws.Range("R8").Value2 = 1.946127481055
$ws.Range("S8").Value2 = 0.009951808920583435
$ws.Range("T8").Value2 = 0.009951808920583437
$ws.Range("G9").Value2 = 3.233763666666667
$ws.Range("H9").Value2 = 9.701291000000001
$ws.Range("I9").Value2 = 0.1787554827329337
$ws.Range("J9").Value2 = 0.1787554827329337
$ws.Range("O9").Value2 = 0.8610905203551533
$ws.Range("P9").Value2 = 0.8610905203551533
$ws.Range("R9").Value2 = 30.100758260578
$ws.Range("S9").Value2 = 0.1539246516428385
$ws.Range("T9").Value2 = 0.1539246516428385
$ws.Range("G10").Value2 = 3.233763666666667
$ws.Range("H10").Value2 = 9.701291000000001
$ws.Range("I10").Value2 = 0.1787554827329337
$ws.Range("J10").Value2 = 0.1787554827329337
$ws.Range("O10").Value2 = 0.08323673177477579
$ws.Range("P10").Value2 = 0.08323673177477578
$ws.Range("Q10").Value2 = 0.3232966004962223
$ws.Range("R10").Value2 = 2.909669404466001
$ws.Range("S10").Value2 = 0.01487902216951176
$ws.Range("T10").Value2 = 0.01487902216951177
